# Adding bioassay 2 PAM results
# - New worksheet "PAM" at the end of the workbook with Number/Date/Station/Other/Notes/FvFm columns
# - "All Nutrients zeroes": clear the leftover highlight-style from the U column DIV/0 + PAM-era rows
#   (rows 6-9, 14-17, 22-25) so they fall back to the default (unstyled) cell format
# - Make "PAM" the active / selected sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Tidy up "All Nutrients zeroes" (sheet7 / U column) - these cells lose
#    their now-stale fill/format (s="14") and fall back to the default style.
# ---------------------------------------------------------------------------
$zeroesSheet = $wb.Worksheets.Item("All Nutrients zeroes")
$clearRows = 6,7,8,9,14,15,16,17,22,23,24,25
foreach ($r in $clearRows) {
    $zeroesSheet.Cells.Item($r, 21).ClearFormats()
}

# ---------------------------------------------------------------------------
# 2. Add the new "PAM" sheet as the last tab in the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$pam = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$pam.Name = "PAM"

# Header row
$pam.Range("A1").Value = "Number"
$pam.Range("B1").Value = "Date"
$pam.Range("C1").Value = "Station"
$pam.Range("D1").Value = "Other"
$pam.Range("E1").Value = "Notes"
$pam.Range("F1").Value = "FvFm"

# Data rows: Number, Date(serial), Station, Other/treatment, FvFm
$data = @(
    @(1, 45083, "Clambank Landing", "T_0", 0.621),
    @(2, 45083, "Clambank Landing", "T_0", 0.585),
    @(3, 45083, "Clambank Landing", "T_0", 0.635),
    @(4, 45083, "Clambank Landing", "T_0", 0.578),
    @(5, 45083, "Clambank Landing", "T_0", 0.607),
    @(6, 45085, "Clambank Landing", "Control", 0.658),
    @(7, 45085, "Clambank Landing", "Control", 0.545),
    @(8, 45085, "Clambank Landing", "Control", 0.73),
    @(9, 45085, "Clambank Landing", "Control", 0.574),
    @(10, 45085, "Clambank Landing", "Control", 0.564),
    @(11, 45085, "Clambank Landing", "DIN", 0.664),
    @(12, 45085, "Clambank Landing", "DIN", 0.686),
    @(13, 45085, "Clambank Landing", "DIN", 0.673),
    @(14, 45085, "Clambank Landing", "DIN", 0.667),
    @(15, 45085, "Clambank Landing", "DIN", 0.655),
    @(16, 45085, "Clambank Landing", "LP", 0.585),
    @(17, 45085, "Clambank Landing", "LP", 0.582),
    @(18, 45085, "Clambank Landing", "LP", 0.641),
    @(19, 45085, "Clambank Landing", "LP", 0.682),
    @(20, 45085, "Clambank Landing", "LP", 0.578),
    @(21, 45085, "Clambank Landing", "HP", 0.617),
    @(22, 45085, "Clambank Landing", "HP", 0.52),
    @(23, 45085, "Clambank Landing", "HP", 0.726),
    @(24, 45085, "Clambank Landing", "HP", 0.594),
    @(25, 45085, "Clambank Landing", "HP", 0.574),
    @(26, 45085, "Clambank Landing", "DIN_LP", 0.678),
    @(27, 45085, "Clambank Landing", "DIN_LP", 0.727),
    @(28, 45085, "Clambank Landing", "DIN_LP", 0.752),
    @(29, 45085, "Clambank Landing", "DIN_LP", 0.849),
    @(30, 45085, "Clambank Landing", "DIN_LP", 0.928),
    @(31, 45085, "Clambank Landing", "DIN_HP", 0.689),
    @(32, 45085, "Clambank Landing", "DIN_HP", 0.692),
    @(33, 45085, "Clambank Landing", "DIN_HP", 0.705),
    @(34, 45085, "Clambank Landing", "DIN_HP", 0.69),
    @(35, 45085, "Clambank Landing", "DIN_HP", 0.773)
)

$row = 2
foreach ($entry in $data) {
    $pam.Cells.Item($row, 1).Value = $entry[0]
    $bcell = $pam.Cells.Item($row, 2)
    $bcell.Value = $entry[1]
    $bcell.NumberFormat = "d-mmm-yy"
    $pam.Cells.Item($row, 3).Value = $entry[2]
    $pam.Cells.Item($row, 4).Value = $entry[3]
    $pam.Cells.Item($row, 6).Value = $entry[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3. Make PAM the active sheet / active cell, matching the saved view state.
# ---------------------------------------------------------------------------
$pam.Activate()
$pam.Range("F37").Select()
